$wb = $excel.ActiveWorkbook

# Sheet 1: Pediatric Vaccine - fix embedded newline in packaging text
$ws1 = $wb.Worksheets.Item("Pediatric Vaccine ")
$ws1.Range("D8").Value = "5 pack - 1 dose T-L syringes. No Needle"

# Sheet 2: Adult Vaccine - fix embedded newline in vaccine name text
$ws2 = $wb.Worksheets.Item("Adult Vaccine ")
$ws2.Range("B14").Value = "Tetanus  Diphtheria Toxoids Adsorbed for Adults No Preservative"

# Sheet 3: Pediatric Influenza Vaccine - fix embedded newlines
$ws3 = $wb.Worksheets.Item("Pediatric Influenza Vaccine ")
$ws3.Range("B3").Value = "Fluzone Pediatric dose No Preservative"
$ws3.Range("B6").Value = "Fluarix Preservative-Free"
$ws3.Range("B9").Value = "FluMist No Preservative"
$ws3.Range("B10").Value = "Afluria No Preservative"
$ws3.Range("H10").Value = "Merck (CSL product)"
$ws3.Range("H11").Value = "Merck (CSL product)"
$ws3.Range("B12").Value = "Afluria No Preservative"
$ws3.Range("H12").Value = "Merck (CSL product)"

# Sheet 4: Adult Influenza Vaccine - fix embedded newlines
$ws4 = $wb.Worksheets.Item("Adult Influenza Vaccine ")
$ws4.Range("B5").Value = "Agriflu No Preservative"
$ws4.Range("B7").Value = "Fluvirin Preservative-free"
$ws4.Range("B8").Value = "Fluraix Preservative-free"
$ws4.Range("B10").Value = "Flumist No Preservative"
